# Updated symbol list on Wed Dec 21 18:41:26 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) cells: numeric-looking strings must stay text ---
# Force text number format first so Excel doesn't coerce these into numbers.
$priceCells = @(
    "D2", "D4", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13",
    "D14", "D15", "D16", "D20", "D21", "D22", "D23", "D25", "D26",
    "D40", "D41", "D42", "D43", "D44", "D45", "D48"
)
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value  = "247.79"
$ws.Range("D4").Value  = "5.235"
$ws.Range("D6").Value  = "3.418"
$ws.Range("D7").Value  = "6.313"
$ws.Range("D8").Value  = "0.8071"
$ws.Range("D9").Value  = "0.8659"
$ws.Range("D10").Value = "0.1416"
$ws.Range("D11").Value = "0.07433"
$ws.Range("D12").Value = "0.03050"
$ws.Range("D13").Value = "0.03077"
$ws.Range("D14").Value = "0.09395"
$ws.Range("D15").Value = "3.877"
$ws.Range("D16").Value = "0.001577"
$ws.Range("D20").Value = "0.006438"
$ws.Range("D21").Value = "0.005040"
$ws.Range("D22").Value = "0.0009961"
$ws.Range("D23").Value = "0.0001500"
$ws.Range("D25").Value = "2.199"
$ws.Range("D26").Value = "0.3246"
$ws.Range("D40").Value = "0.03956"
$ws.Range("D41").Value = "0.1067"
$ws.Range("D42").Value = "0.002731"
$ws.Range("D43").Value = "0.003044"
$ws.Range("D44").Value = "0.007970"
$ws.Range("D45").Value = "0.00005590"
$ws.Range("D48").Value = "0.2023"

# --- Row 19: "Worstin24h" suffix removed from E19 ---
$ws.Range("E19").Value = "18OneONE"

# --- Rows 41-43: coin rows shuffled (KickToken -> BKEXToken -> CEJI -> KickToken) ---
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"
